$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.354, -0.288, 0.008, -102, 0, 63),
    @(-0.301, 0.108, -0.086, -50, 27, -61),
    @(-0.273, -0.268, 0.645, 121, -70, -130),
    @(0.252, -0.29, -0.139, -30, -23, -82),
    @(0.237, 0.068, -0.126, 10, -33, -162),
    @(-0.014, 0.154, -0.254, -55, 31, -1),
    @(-0.221, 0.05, -0.333, -17, -80, -30),
    @(0.303, -0.153, 0.188, -81, -85, -29),
    @(0.384, 0.194, 0.316, -86, -41, -126),
    @(-0.131, 0.18, -0.305, -114, 22, 27)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
